$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix calculator choice logic ---------------------------------------
# The "Liczba Zgadnięć" (number of guesses) results table used to start
# at row 16 (header) / row 17 (first data row), leaving a big empty gap
# (rows 7-15) below the scenario-setup block above it. Close that gap so
# the table starts right after the setup block: header on row 9, 20 data
# rows on rows 10-29.
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart()

# Remember the chart's current on-sheet position/size (absolute points)
# before we shift rows, so we can re-anchor it at the same visual spot
# relative to the row that used to be 37 (now 30) after the rows above
# it disappear.
$chartTop0 = $co.Top()

# Delete the 7 superfluous blank rows (7:13) so the table below collapses
# upward by 7 rows (16 -> 9, 17:36 -> 10:29).
$ws.Rows("7:13").Delete()

# Re-anchor the chart so its top-left corner keeps the same position
# relative to the rows that are still above it (shift up by the height
# of the 7 removed rows).
$rowsRemoved = 7
$rowHeightPts = $ws.Rows(1).RowHeight()
$co.Top = $chartTop0 - ($rowsRemoved * $rowHeightPts)

# Update each plotted series' source formula to point at the table's new
# location (columns B..I, header row 9, data rows 10:29).
$cols = @("B", "C", "D", "E", "F", "G", "H", "I")
$sc = $chart.SeriesCollection()
for ($i = 1; $i -le $cols.Length; $i++) {
    $col = $cols[$i - 1]
    $ser = $sc.Item($i)
    $ser.Formula = "=SERIES(Arkusz1!`$" + $col + "`$9,,Arkusz1!`$" + $col + "`$10:`$" + $col + "`$29," + $i + ")"
}

# Move the active selection to match the post-edit cursor position.
$ws.Range("N25").Select()
